$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Level" column (C) for great-grandparent, grandmother, and family rows
# from "nothing" to "hard"
$ws.Range("C2").Value = "hard"
$ws.Range("C3").Value = "hard"
$ws.Range("C6").Value = "hard"

# Move selection to C6
$ws.Range("C6").Select()
